# version 2 - server mode added, both working.
#
# This script reproduces the authoring changes made to galilean.xlsx:
#  - Tuesday sheet: fill in C9 ("anna") / D9 ("steve"), move the selection
#    to D9 (it is no longer the active tab).
#  - pictures sheet: add a value (2000) to B3, and leave its selection on
#    B13 (it briefly becomes active while we touch it, but is not the
#    sheet that ends up selected when the workbook is saved).
#  - Thursday sheet: fill in several more name cells, add a brand new
#    row 16, and make Thursday the active/selected tab with C20 selected.

$wb = $excel.ActiveWorkbook

# ---- Tuesday (was the active tab, loses that status) ----------------------
$wsTuesday = $wb.Worksheets.Item("Tuesday")
$wsTuesday.Range("C9").Value = "anna"
$wsTuesday.Range("D9").Value = "steve"
$wsTuesday.Range("D9").Select()

# ---- pictures ---------------------------------------------------------
$wsPictures = $wb.Worksheets.Item("pictures")
$wsPictures.Range("B3").Value = 2000
$wsPictures.Range("B13").Select()

# ---- Thursday (becomes the new active tab) ---------------------------
$wsThursday = $wb.Worksheets.Item("Thursday")
$wsThursday.Range("C9").Value = "anna"
$wsThursday.Range("D9").Value = "steve"
$wsThursday.Range("B13").Value = "anna"
$wsThursday.Range("D13").Value = "steve"
$wsThursday.Range("C14").Value = "jack"
$wsThursday.Range("A15").Value = 21
$wsThursday.Range("B15").Value = "anna"
$wsThursday.Range("D15").Value = "anna"
$wsThursday.Range("A16").Value = 22
$wsThursday.Range("B16").Value = "anna"
$wsThursday.Range("C16").Value = "anna"

$wsThursday.Activate()
$wsThursday.Range("C20").Select()
